$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.457486347212352
$ws.Range("C2").Value = 0.1444092747045715
$ws.Range("D2").Value = 0.0761708349506165
$ws.Range("E2").Value = 0.03034151283318853
$ws.Range("G2").Value = 0.002622539093792955
$ws.Range("I2").Value = 2.550472602129958
$ws.Range("K2").Value = 1.416355803092387
$ws.Range("L2").Value = 0.2380638193812388
$ws.Range("M2").Value = 0.3287198858151115
$ws.Range("N2").Value = 4.303466944968847
$ws.Range("B3").Value = 1.419144728633569
$ws.Range("C3").Value = 0.1330234071336918
$ws.Range("D3").Value = 0.06945171253869376
$ws.Range("E3").Value = 0.03040163909198768
$ws.Range("G3").Value = 0.002627427494280126
$ws.Range("I3").Value = 2.533809843593133
$ws.Range("K3").Value = 1.36688427237155
$ws.Range("L3").Value = 0.2353324112451105
$ws.Range("M3").Value = 0.3218263764707565
$ws.Range("N3").Value = 4.291761579579301
$ws.Range("B4").Value = 1.396505353097922
$ws.Range("C4").Value = 0.1261117126981617
$ws.Range("D4").Value = 0.06536671740104794
$ws.Range("E4").Value = 0.03044274976307659
$ws.Range("G4").Value = 0.002630587460974521
$ws.Range("I4").Value = 2.524388391056462
$ws.Range("K4").Value = 1.337402948019502
$ws.Range("L4").Value = 0.2337702712474226
$ws.Range("M4").Value = 0.3177845164548003
$ws.Range("N4").Value = 4.285322170794245
$ws.Range("B5").Value = 1.387506299295325
$ws.Range("C5").Value = 0.123314831707205
$ws.Range("D5").Value = 0.06371213979205947
$ws.Range("E5").Value = 0.03046055942131409
$ws.Range("G5").Value = 0.002631915154110536
$ws.Range("I5").Value = 2.520752427417818
$ws.Range("K5").Value = 1.325613394670199
$ws.Range("L5").Value = 0.2331625906580115
$ws.Range("M5").Value = 0.3161853604115628
$ws.Range("N5").Value = 4.282885761800543
$ws.Range("B6").Value = 1.386025696197123
$ws.Range("C6").Value = 0.1228515939450858
$ws.Range("D6").Value = 0.06343800468570748
$ws.Range("E6").Value = 0.03046358059629384
$ws.Range("G6").Value = 0.002632138035315136
$ws.Range("I6").Value = 2.520160954494116
$ws.Range("K6").Value = 1.323669280477787
$ws.Range("L6").Value = 0.2330634314571256
$ws.Range("M6").Value = 0.3159227158218627
$ws.Range("N6").Value = 4.282492525419869
$ws.Range("B7").Value = 1.396383071269412
$ws.Range("C7").Value = 0.1260739135743734
$ws.Range("D7").Value = 0.06534436247537201
$ws.Range("E7").Value = 0.0304429856683055
$ws.Range("G7").Value = 0.002630605204568255
$ws.Range("I7").Value = 2.524338532099691
$ws.Range("K7").Value = 1.33724304245186
$ws.Range("L7").Value = 0.2337619588246724
$ws.Range("M7").Value = 0.3177627556419793
$ws.Range("N7").Value = 4.285288552987851
$ws.Range("B8").Value = 1.444078652267592
$ws.Range("C8").Value = 0.1404668234505948
$ws.Range("D8").Value = 0.07384559119533662
$ws.Range("E8").Value = 0.030361375540916
$ws.Range("G8").Value = 0.002624191804137962
$ws.Range("I8").Value = 2.544559122056711
$ws.Range("K8").Value = 1.399112017461135
$ws.Range("L8").Value = 0.2370981676364536
$ws.Range("M8").Value = 0.3263033796569559
$ws.Range("N8").Value = 4.299275541164718
$ws.Range("B9").Value = 1.544790525173596
$ws.Range("C9").Value = 0.1693328293615934
$ws.Range("D9").Value = 0.09084433753359633
$ws.Range("E9").Value = 0.03023450723004872
$ws.Range("G9").Value = 0.002612866450542865
$ws.Range("I9").Value = 2.590648948768631
$ws.Range("K9").Value = 1.527568221164387
$ws.Range("L9").Value = 0.2445533730668075
$ws.Range("M9").Value = 0.344568477994379
$ws.Range("N9").Value = 4.332654761678896
$ws.Range("B10").Value = 1.623199401975853
$ws.Range("C10").Value = 0.1909514162072696
$ws.Range("D10").Value = 0.1035427139092775
$ws.Range("E10").Value = 0.03016139071773161
$ws.Range("G10").Value = 0.002605300029615422
$ws.Range("I10").Value = 2.628460987120661
$ws.Range("K10").Value = 1.626353006579393
$ws.Range("L10").Value = 0.2505893241053059
$ws.Range("M10").Value = 0.3589190935674509
$ws.Range("N10").Value = 4.360835820295534
$ws.Range("B11").Value = 1.659837567164971
$ws.Range("C11").Value = 0.2008798557406237
$ws.Range("D11").Value = 0.1093671465495731
$ws.Range("E11").Value = 0.03013246582046891
$ws.Range("G11").Value = 0.002602019845882271
$ws.Range("I11").Value = 2.646526339628309
$ws.Range("K11").Value = 1.672263948128204
$ws.Range("L11").Value = 0.2534570596076549
$ws.Range("M11").Value = 0.3656513301165347
$ws.Range("N11").Value = 4.374457050637602
$ws.Range("B12").Value = 1.673851451905307
$ws.Range("C12").Value = 0.2046533323395465
$ws.Range("D12").Value = 0.1115797393034086
$ws.Range("E12").Value = 0.0301221342131992
$ws.Range("G12").Value = 0.002600800856654246
$ws.Range("I12").Value = 2.653491908810693
$ws.Range("K12").Value = 1.68979008448531
$ws.Range("L12").Value = 0.2545605585720239
$ws.Range("M12").Value = 0.3682300874535613
$ws.Range("N12").Value = 4.379730782527275
$ws.Range("B13").Value = 1.670827086895372
$ws.Range("C13").Value = 0.2038400289049775
$ws.Range("D13").Value = 0.111102904295592
$ws.Range("E13").Value = 0.03012433169273354
$ws.Range("G13").Value = 0.002601062360275237
$ws.Range("I13").Value = 2.651986201771365
$ws.Range("K13").Value = 1.686009250146242
$ws.Range("L13").Value = 0.2543221195701477
$ws.Range("M13").Value = 0.3676733971140109
$ws.Range("N13").Value = 4.378589839058407
$ws.Range("B14").Value = 1.660987693718255
$ws.Range("C14").Value = 0.201190023962539
$ws.Range("D14").Value = 0.1095490367635392
$ws.Range("E14").Value = 0.03013160338542331
$ws.Range("G14").Value = 0.00260191909585465
$ws.Range("I14").Value = 2.647096902016528
$ws.Range("K14").Value = 1.673703010785459
$ws.Range("L14").Value = 0.2535474933274173
$ws.Range("M14").Value = 0.3658628963495829
$ws.Range("N14").Value = 4.374888603216675
$ws.Range("B15").Value = 1.654978999859452
$ws.Range("C15").Value = 0.1995686236073766
$ws.Range("D15").Value = 0.108598164125624
$ws.Range("E15").Value = 0.0301361384052381
$ws.Range("G15").Value = 0.002602446881357935
$ws.Range("I15").Value = 2.644118302072911
$ws.Range("K15").Value = 1.666183423089478
$ws.Range("L15").Value = 0.2530752983014821
$ws.Range("M15").Value = 0.3647577436536906
$ws.Range("N15").Value = 4.372636564985044
$ws.Range("B16").Value = 1.620824554468527
$ws.Range("C16").Value = 0.1903044845678608
$ws.Range("D16").Value = 0.1031630473383558
$ws.Range("E16").Value = 0.03016336811770737
$ws.Range("G16").Value = 0.002605517642984498
$ws.Range("I16").Value = 2.627297794832074
$ws.Range("K16").Value = 1.623372262908134
$ws.Range("L16").Value = 0.25040436592775
$ws.Range("M16").Value = 0.3584832375914146
$ws.Range("N16").Value = 4.359961804416798
$ws.Range("B17").Value = 1.60012050797269
$ws.Range("C17").Value = 0.1846455111441117
$ws.Range("D17").Value = 0.09984112806444045
$ws.Range("E17").Value = 0.03018118188145147
$ws.Range("G17").Value = 0.002607442812847975
$ws.Range("I17").Value = 2.617200581958627
$ws.Range("K17").Value = 1.597358773057749
$ws.Range("L17").Value = 0.2487970791449072
$ws.Range("M17").Value = 0.3546863407026351
$ws.Range("N17").Value = 4.352391828143084
$ws.Range("B18").Value = 1.588303325183119
$ws.Range("C18").Value = 0.1813994540922579
$ws.Range("D18").Value = 0.09793494705736805
$ws.Range("E18").Value = 0.03019183613843301
$ws.Range("G18").Value = 0.002608565358099805
$ws.Range("I18").Value = 2.611474281241669
$ws.Range("K18").Value = 1.582488065101245
$ws.Range("L18").Value = 0.2478840880425963
$ws.Range("M18").Value = 0.3525216702554275
$ws.Range("N18").Value = 4.348113196507427
$ws.Range("B19").Value = 1.584317888744579
$ws.Range("C19").Value = 0.1803019041643665
$ws.Range("D19").Value = 0.09729031595773563
$ws.Range("E19").Value = 0.03019551366373463
$ws.Range("G19").Value = 0.00260894805357053
$ws.Range("I19").Value = 2.609549416885415
$ws.Range("K19").Value = 1.57746880369794
$ws.Range("L19").Value = 0.2475769358688638
$ws.Range("M19").Value = 0.3517920465760795
$ws.Range("N19").Value = 4.34667746618004
$ws.Range("B20").Value = 1.602315043649071
$ws.Range("C20").Value = 0.1852470022115824
$ws.Range("D20").Value = 0.1001942855205016
$ws.Range("E20").Value = 0.03017924333705801
$ws.Range("G20").Value = 0.002607236298871733
$ws.Range("I20").Value = 2.618267026055733
$ws.Range("K20").Value = 1.600118471934081
$ws.Range("L20").Value = 0.2489669895710165
$ws.Range("M20").Value = 0.3550885388666245
$ws.Range("N20").Value = 4.353189856630337
$ws.Range("B21").Value = 1.663873964568836
$ws.Range("C21").Value = 0.2019680181742558
$ws.Range("D21").Value = 0.1100052541914778
$ws.Range("E21").Value = 0.03012945065738037
$ws.Range("G21").Value = 0.002601666824668447
$ws.Range("I21").Value = 2.64852962347328
$ws.Range("K21").Value = 1.677313828425895
$ws.Range("L21").Value = 0.2537745434240293
$ws.Range("M21").Value = 0.3663938857603526
$ws.Range("N21").Value = 4.375972603419626
$ws.Range("B22").Value = 1.704921366211636
$ws.Range("C22").Value = 0.2129766866646605
$ws.Range("D22").Value = 0.1164581861994236
$ws.Range("E22").Value = 0.03010053071670171
$ws.Range("G22").Value = 0.002598161701575507
$ws.Range("I22").Value = 2.669034462715189
$ws.Range("K22").Value = 1.728585665596881
$ws.Range("L22").Value = 0.2570188583203219
$ws.Range("M22").Value = 0.3739540109763126
$ws.Range("N22").Value = 4.391536839597507
$ws.Range("B23").Value = 1.682938884643647
$ws.Range("C23").Value = 0.2070936957429126
$ws.Range("D23").Value = 0.1130103522961576
$ws.Range("E23").Value = 0.03011563499427528
$ws.Range("G23").Value = 0.002600020153612118
$ws.Range("I23").Value = 2.658024068535425
$ws.Range("K23").Value = 1.701145634033026
$ws.Range("L23").Value = 0.2552779411680035
$ws.Range("M23").Value = 0.3699033235718403
$ws.Range("N23").Value = 4.383168069659746
$ws.Range("B24").Value = 1.601322626689921
$ws.Range("C24").Value = 0.1849750451900718
$ws.Range("D24").Value = 0.1000346117239985
$ws.Range("E24").Value = 0.03018011846689239
$ws.Range("G24").Value = 0.002607329614810707
$ws.Range("I24").Value = 2.617784641483595
$ws.Range("K24").Value = 1.598870548007341
$ws.Range("L24").Value = 0.248890138620169
$ws.Range("M24").Value = 0.3549066483408509
$ws.Range("N24").Value = 4.352828839174833
$ws.Range("B25").Value = 1.516772626542149
$ws.Range("C25").Value = 0.1614531095615632
$ws.Range("D25").Value = 0.08620964634913264
$ws.Range("E25").Value = 0.03026529086282626
$ws.Range("G25").Value = 0.002615797180367344
$ws.Range("I25").Value = 2.577488666917858
$ws.Range("K25").Value = 1.492047752186068
$ws.Range("L25").Value = 0.2424385987592359
$ws.Range("M25").Value = 0.3394642244096531
$ws.Range("N25").Value = 4.322984692065376
